# Apply weekly update: insert a new price record as row 17, shifting
# the existing rows 17-47 down to 18-48 (dimension grows from T47 to T48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17; this pushes old rows 17..47 down to 18..48
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record's data.
# Columns A-C, E-L are the same as the (old, now-shifted) row 17/new row 18
# record for Castle Brite / Primera, Region Metropolitana market entry.
$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C17").Value = "Metropolitana"
$ws.Range("D17").Value = 44525
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100103
$ws.Range("H17").Value = "Frutos de hueso (carozo)"
$ws.Range("I17").Value = 100103003
$ws.Range("J17").Value = "Damasco"
$ws.Range("K17").Value = "Castle Brite"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 12
$ws.Range("N17").Value = 580000
$ws.Range("O17").Value = 580000
$ws.Range("P17").Value = 580000
$ws.Range("Q17").Value = "`$/bins (500 kilos)"
$ws.Range("R17").Value = "Paine"
$ws.Range("S17").Value = 1160
$ws.Range("T17").Value = 500
